# Updates the cryptos list (Price and Volume(1h) columns) on sheet1,
# matching the latest GitHub Actions scrape refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.474.98"; E = "  +2.13%  " },
    @{ Row = 3; D = "1.839.73"; E = "  +1.46%  " },
    @{ Row = 4; D = $null; E = "  +1.22%  " },
    @{ Row = 5; D = "314.70"; E = "  +1.71%  " },
    @{ Row = 6; D = $null; E = "  +1.01%  " },
    @{ Row = 7; D = "0.4748"; E = "  +1.74%  " },
    @{ Row = 8; D = "0.3694"; E = "  +1.04%  " },
    @{ Row = 9; D = "0.07467"; E = "  +1.56%  " },
    @{ Row = 10; D = "0.8857"; E = "  +1.97%  " },
    @{ Row = 11; D = "20.46"; E = "  +0.73%  " },
    @{ Row = 12; D = "1.877.43"; E = "  +4.45%  " },
    @{ Row = 13; D = "0.07350"; E = "  +3.69%  " },
    @{ Row = 14; D = "5.449"; E = "  +1.19%  " },
    @{ Row = 15; D = "93.33"; E = "  +1.87%  " },
    @{ Row = 16; D = "6.586"; E = "  +1.15%  " },
    @{ Row = 17; D = $null; E = "  +1.04%  " },
    @{ Row = 18; D = "0.000008819"; E = "  +1.31%  " },
    @{ Row = 19; D = "1.012"; E = "  +1.04%  " },
    @{ Row = 20; D = $null; E = "  +1.21%  " },
    @{ Row = 21; D = "27.495.98"; E = "  +2.16%  " },
    @{ Row = 22; D = "5.322"; E = "  +0.46%  " },
    @{ Row = 23; D = "10.71"; E = "  +0.72%  " },
    @{ Row = 24; D = "2.092.23"; E = "  +2.87%  " },
    @{ Row = 25; D = "1.911"; E = "  +0.87%  " },
    @{ Row = 26; D = "152.36"; E = "  +1.07%  " },
    @{ Row = 27; D = "18.65"; E = "  +1.74%  " },
    @{ Row = 28; D = $null; E = "  +0.41%  " },
    @{ Row = 29; D = "5.261"; E = "  -0.11%  " },
    @{ Row = 30; D = "117.98"; E = "  +2.44%  " },
    @{ Row = 31; D = "0.09005"; E = "  +0.83%  " },
    @{ Row = 32; D = "0.7564"; E = "  +0.35%  " },
    @{ Row = 33; D = $null; E = "  +2.33%  " },
    @{ Row = 34; D = "4.565"; E = "  +1.69%  " },
    @{ Row = 35; D = "2.954"; E = "  +1.41%  " },
    @{ Row = 36; D = $null; E = "  +1.18%  " },
    @{ Row = 37; D = "1.108"; E = "  +2.23%  " },
    @{ Row = 38; D = "0.05336"; E = "  +1.05%  " },
    @{ Row = 39; D = "0.01958"; E = "  +0.55%  " },
    @{ Row = 40; D = "3.000"; E = "  +0.86%  " },
    @{ Row = 41; D = "7.333"; E = "  +1.72%  " },
    @{ Row = 42; D = "2.398"; E = "  +5.14%  " },
    @{ Row = 43; D = "0.5340"; E = "  +0.76%  " },
    @{ Row = 44; D = $null; E = "  +0.53%  " },
    @{ Row = 45; D = "8.526"; E = "  +1.49%  " },
    @{ Row = 46; D = "0.4920"; E = "  +0.98%  " },
    @{ Row = 47; D = "10.58"; E = "  +1.36%  " },
    @{ Row = 48; D = "1.013"; E = "  +1.15%  " },
    @{ Row = 49; D = "104.91"; E = "  +1.80%  " },
    @{ Row = 50; D = "1.679"; E = "  +1.33%  " },
    @{ Row = 51; D = "0.06320"; E = "  +0.44%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
